# Refresh the cryptos table (Price / Volume(1h) columns) with latest values.
# D-column price cells are stored as text in the source sheet (mixed
# thousand-separator formatting like "37.826.18"), so force NumberFormat
# to "@" before writing so Excel does not reinterpret them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.826.18'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.077.92'
$ws.Range("E3").Value = '  -1.23%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.61'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.31'
$ws.Range("E7").Value = '  +2.59%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +0.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0790'
$ws.Range("E10").Value = '  +1.27%  '
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.81'
$ws.Range("E12").Value = '  +1.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.20'
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.775'
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.36'
$ws.Range("E15").Value = '  +1.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.069.00'
$ws.Range("E16").Value = '  -1.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '37.721.54'
$ws.Range("E17").Value = '  -0.23%  '
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.60'
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0856'
$ws.Range("E20").Value = '  +3.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '228.33'
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("E23").Value = '  +1.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.37'
$ws.Range("E24").Value = '  -1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.64'
$ws.Range("E25").Value = '  +1.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.19'
$ws.Range("E26").Value = '  +2.36%  '
$ws.Range("E27").Value = '  -2.91%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.48'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("E30").Value = '  +1.81%  '
$ws.Range("E31").Value = '  +1.32%  '
$ws.Range("E32").Value = '  +2.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0633'
$ws.Range("E33").Value = '  +0.50%  '
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("E36").Value = '  -1.06%  '
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0984'
$ws.Range("E39").Value = '  -0.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '99.36'
$ws.Range("E40").Value = '  +1.79%  '
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("E42").Value = '  -1.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.73'
$ws.Range("E43").Value = '  +6.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.445.02'
$ws.Range("E44").Value = '  -0.83%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.15'
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.14'
$ws.Range("E46").Value = '  +1.95%  '
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.40'
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.02'
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.269.95'
$ws.Range("E50").Value = '  -1.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '46.73'
$ws.Range("E51").Value = '  +0.65%  '
